$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.039822463303314
$ws.Cells.Item(2, 4).Value = 1.042121712432953
$ws.Cells.Item(2, 5).Value = 1.047808065905107
$ws.Cells.Item(2, 6).Value = 1.057498777878844
$ws.Cells.Item(2, 9).Value = 1.042918224583322
$ws.Cells.Item(2, 10).Value = 1.044912398873431
$ws.Cells.Item(2, 11).Value = 1.044899293951856
$ws.Cells.Item(2, 12).Value = 1.050569669777072
$ws.Cells.Item(2, 13).Value = 1.060233581550295

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.040682827700167
$ws.Cells.Item(3, 4).Value = 1.042759306656269
$ws.Cells.Item(3, 5).Value = 1.048585319699209
$ws.Cells.Item(3, 6).Value = 1.058386113757657
$ws.Cells.Item(3, 9).Value = 1.043150409370334
$ws.Cells.Item(3, 10).Value = 1.045418573047115
$ws.Cells.Item(3, 11).Value = 1.045348109437201
$ws.Cells.Item(3, 12).Value = 1.051158932055789
$ws.Cells.Item(3, 13).Value = 1.06093457893126

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.041240068048999
$ws.Cells.Item(4, 4).Value = 1.043172236401728
$ws.Cells.Item(4, 5).Value = 1.049089122854325
$ws.Cells.Item(4, 6).Value = 1.058961317616562
$ws.Cells.Item(4, 9).Value = 1.043299591956947
$ws.Cells.Item(4, 10).Value = 1.045745966312441
$ws.Cells.Item(4, 11).Value = 1.045638197542667
$ws.Cells.Item(4, 12).Value = 1.051540431010173
$ws.Cells.Item(4, 13).Value = 1.061388569452126

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.04147445656412
$ws.Cells.Item(5, 4).Value = 1.04334591746405
$ws.Cells.Item(5, 5).Value = 1.049301127768389
$ws.Cells.Item(5, 6).Value = 1.059203380062493
$ws.Cells.Item(5, 9).Value = 1.043362054614016
$ws.Cells.Item(5, 10).Value = 1.045883568844992
$ws.Cells.Item(5, 11).Value = 1.045760071459705
$ws.Cells.Item(5, 12).Value = 1.05170086132877
$ws.Cells.Item(5, 13).Value = 1.061579520888876

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.041513818694332
$ws.Cells.Item(6, 4).Value = 1.043375084238061
$ws.Cells.Item(6, 5).Value = 1.04933673635999
$ws.Cells.Item(6, 6).Value = 1.059244037828963
$ws.Cells.Item(6, 9).Value = 1.043372527469259
$ws.Cells.Item(6, 10).Value = 1.045906670910369
$ws.Cells.Item(6, 11).Value = 1.045780529952234
$ws.Cells.Item(6, 12).Value = 1.051727801069559
$ws.Cells.Item(6, 13).Value = 1.061611587932272

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041243199471811
$ws.Cells.Item(7, 4).Value = 1.043174556803119
$ws.Cells.Item(7, 5).Value = 1.049091954866828
$ws.Cells.Item(7, 6).Value = 1.058964551099543
$ws.Cells.Item(7, 9).Value = 1.04330042758387
$ws.Cells.Item(7, 10).Value = 1.045747805097744
$ws.Cells.Item(7, 11).Value = 1.045639826340302
$ws.Cells.Item(7, 12).Value = 1.051542574500424
$ws.Cells.Item(7, 13).Value = 1.061391120587847

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040113117443816
$ws.Cells.Item(8, 4).Value = 1.042337114092236
$ws.Cells.Item(8, 5).Value = 1.048070561704961
$ws.Cells.Item(8, 6).Value = 1.057798441139049
$ws.Cells.Item(8, 9).Value = 1.042996910807385
$ws.Cells.Item(8, 10).Value = 1.045083490085615
$ws.Cells.Item(8, 11).Value = 1.045051039932464
$ws.Cells.Item(8, 12).Value = 1.050768770071914
$ws.Cells.Item(8, 13).Value = 1.060470403490928

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.038125871715058
$ws.Cells.Item(9, 4).Value = 1.040864291502417
$ws.Cells.Item(9, 5).Value = 1.046277460487591
$ws.Cells.Item(9, 6).Value = 1.055751637725469
$ws.Cells.Item(9, 9).Value = 1.042454018258113
$ws.Cells.Item(9, 10).Value = 1.04391189969904
$ws.Cells.Item(9, 11).Value = 1.044011082560252
$ws.Cells.Item(9, 12).Value = 1.049406872369339
$ws.Cells.Item(9, 13).Value = 1.058851101320852

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.036803889301303
$ws.Cells.Item(10, 4).Value = 1.039884430488777
$ws.Cells.Item(10, 5).Value = 1.045086683406827
$ws.Cells.Item(10, 6).Value = 1.054392609441817
$ws.Cells.Item(10, 9).Value = 1.042086718776604
$ws.Cells.Item(10, 10).Value = 1.043130250922797
$ws.Cells.Item(10, 11).Value = 1.043316207619356
$ws.Cells.Item(10, 12).Value = 1.048500127417078
$ws.Cells.Item(10, 13).Value = 1.057773755320333

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036232150528977
$ws.Cells.Item(11, 4).Value = 1.039460639974644
$ws.Cells.Item(11, 5).Value = 1.04457218214707
$ws.Cells.Item(11, 6).Value = 1.053805464082878
$ws.Cells.Item(11, 9).Value = 1.04192640975119
$ws.Cells.Item(11, 10).Value = 1.042791663911963
$ws.Cells.Item(11, 11).Value = 1.043014961594752
$ws.Cells.Item(11, 12).Value = 1.048107795411837
$ws.Cells.Item(11, 13).Value = 1.057307791052855

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036019886192581
$ws.Cells.Item(12, 4).Value = 1.039303301472834
$ws.Cells.Item(12, 5).Value = 1.044381242448094
$ws.Cells.Item(12, 6).Value = 1.053587572721888
$ws.Cells.Item(12, 9).Value = 1.041866674407089
$ws.Cells.Item(12, 10).Value = 1.042665879606897
$ws.Cells.Item(12, 11).Value = 1.042903012421055
$ws.Cells.Item(12, 12).Value = 1.047962111456734
$ws.Cells.Item(12, 13).Value = 1.057134792901743

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036065412842439
$ws.Cells.Item(13, 4).Value = 1.039337047642014
$ws.Cells.Item(13, 5).Value = 1.044422191987024
$ws.Cells.Item(13, 6).Value = 1.053634302053058
$ws.Cells.Item(13, 9).Value = 1.041879496401308
$ws.Cells.Item(13, 10).Value = 1.042692861563095
$ws.Cells.Item(13, 11).Value = 1.042927028277552
$ws.Cells.Item(13, 12).Value = 1.047993359073254
$ws.Cells.Item(13, 13).Value = 1.057171897869321

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036214602544176
$ws.Cells.Item(14, 4).Value = 1.039447632761562
$ws.Cells.Item(14, 5).Value = 1.044556395554789
$ws.Cells.Item(14, 6).Value = 1.053787449012932
$ws.Cells.Item(14, 9).Value = 1.041921475869031
$ws.Cells.Item(14, 10).Value = 1.042781266902561
$ws.Cells.Item(14, 11).Value = 1.043005708915187
$ws.Cells.Item(14, 12).Value = 1.048095752192109
$ws.Cells.Item(14, 13).Value = 1.057293489295597

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.036306537188104
$ws.Cells.Item(15, 4).Value = 1.039515778046232
$ws.Cells.Item(15, 5).Value = 1.044639105245606
$ws.Cells.Item(15, 6).Value = 1.053881834541162
$ws.Cells.Item(15, 9).Value = 1.041947315725378
$ws.Cells.Item(15, 10).Value = 1.042835733999071
$ws.Cells.Item(15, 11).Value = 1.043054179667023
$ws.Cells.Item(15, 12).Value = 1.048158846058093
$ws.Cells.Item(15, 13).Value = 1.057368416662018

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03684184834983
$ws.Cells.Item(16, 4).Value = 1.039912566658297
$ws.Cells.Item(16, 5).Value = 1.045120852761581
$ws.Cells.Item(16, 6).Value = 1.054431604372295
$ws.Cells.Item(16, 9).Value = 1.042097331351604
$ws.Cells.Item(16, 10).Value = 1.0431527192352
$ws.Cells.Item(16, 11).Value = 1.043336192830019
$ws.Cells.Item(16, 12).Value = 1.04852617152264
$ws.Cells.Item(16, 13).Value = 1.057804691209919

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.037177820233703
$ws.Cells.Item(17, 4).Value = 1.04016159559059
$ws.Cells.Item(17, 5).Value = 1.04542333948814
$ws.Cells.Item(17, 6).Value = 1.054776815916071
$ws.Cells.Item(17, 9).Value = 1.042191093799892
$ws.Cells.Item(17, 10).Value = 1.04335152229975
$ws.Cells.Item(17, 11).Value = 1.043512996534573
$ws.Cells.Item(17, 12).Value = 1.048756664905437
$ws.Cells.Item(17, 13).Value = 1.058078498676516

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.037373853264066
$ws.Cells.Item(18, 4).Value = 1.040306897675695
$ws.Cells.Item(18, 5).Value = 1.045599882156258
$ws.Cells.Item(18, 6).Value = 1.05497829948313
$ws.Cells.Item(18, 9).Value = 1.042245661655281
$ws.Cells.Item(18, 10).Value = 1.043467468339224
$ws.Cells.Item(18, 11).Value = 1.043616088272853
$ws.Cells.Item(18, 12).Value = 1.0488911361233
$ws.Cells.Item(18, 13).Value = 1.058238257325913

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.03744070666445
$ws.Cells.Item(19, 4).Value = 1.040356449975116
$ws.Cells.Item(19, 5).Value = 1.045660096786794
$ws.Cells.Item(19, 6).Value = 1.055047021745821
$ws.Cells.Item(19, 9).Value = 1.042264247134528
$ws.Cells.Item(19, 10).Value = 1.043507000798235
$ws.Cells.Item(19, 11).Value = 1.043651233940576
$ws.Cells.Item(19, 12).Value = 1.048936992096715
$ws.Cells.Item(19, 13).Value = 1.058292739565779

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037141766760198
$ws.Cells.Item(20, 4).Value = 1.040134872194685
$ws.Cells.Item(20, 5).Value = 1.045390874416166
$ws.Cells.Item(20, 6).Value = 1.054739764781927
$ws.Cells.Item(20, 9).Value = 1.042181046610223
$ws.Cells.Item(20, 10).Value = 1.043330193888797
$ws.Cells.Item(20, 11).Value = 1.043494030768657
$ws.Cells.Item(20, 12).Value = 1.048731932206922
$ws.Cells.Item(20, 13).Value = 1.058049116400302

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.03617066698776
$ws.Cells.Item(21, 4).Value = 1.039415066084122
$ws.Cells.Item(21, 5).Value = 1.044516871270373
$ws.Cells.Item(21, 6).Value = 1.053742345499029
$ws.Cells.Item(21, 9).Value = 1.041909119185325
$ws.Cells.Item(21, 10).Value = 1.042755234221105
$ws.Cells.Item(21, 11).Value = 1.042982540881637
$ws.Cells.Item(21, 12).Value = 1.048065598700356
$ws.Cells.Item(21, 13).Value = 1.057257681382699

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.035560705676821
$ws.Cells.Item(22, 4).Value = 1.038962936907341
$ws.Cells.Item(22, 5).Value = 1.043968329318529
$ws.Cells.Item(22, 6).Value = 1.0531163897361
$ws.Cells.Item(22, 9).Value = 1.041737051988009
$ws.Cells.Item(22, 10).Value = 1.042393630634234
$ws.Cells.Item(22, 11).Value = 1.042660640908984
$ws.Cells.Item(22, 12).Value = 1.047646913000172
$ws.Cells.Item(22, 13).Value = 1.056760548377172

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.035883999612591
$ws.Cells.Item(23, 4).Value = 1.039202576729481
$ws.Cells.Item(23, 5).Value = 1.044259028449818
$ws.Cells.Item(23, 6).Value = 1.053448110039895
$ws.Cells.Item(23, 9).Value = 1.041828371665052
$ws.Cells.Item(23, 10).Value = 1.042585332956637
$ws.Cells.Item(23, 11).Value = 1.042831314781901
$ws.Cells.Item(23, 12).Value = 1.047868840671514
$ws.Cells.Item(23, 13).Value = 1.057024042580241

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037158057568233
$ws.Cells.Item(24, 4).Value = 1.040146947199755
$ws.Cells.Item(24, 5).Value = 1.045405543653493
$ws.Cells.Item(24, 6).Value = 1.054756506200822
$ws.Cells.Item(24, 9).Value = 1.042185586880402
$ws.Cells.Item(24, 10).Value = 1.043339831318166
$ws.Cells.Item(24, 11).Value = 1.043502600689789
$ws.Cells.Item(24, 12).Value = 1.04874310776141
$ws.Cells.Item(24, 13).Value = 1.058062392827586

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.038639126929202
$ws.Cells.Item(25, 4).Value = 1.041244702509481
$ws.Cells.Item(25, 5).Value = 1.046740212697545
$ws.Cells.Item(25, 6).Value = 1.056279822639949
$ws.Cells.Item(25, 9).Value = 1.042595319072201
$ws.Cells.Item(25, 10).Value = 1.044214892902563
$ws.Cells.Item(25, 11).Value = 1.044280218381712
$ws.Cells.Item(25, 12).Value = 1.049758752320429
$ws.Cells.Item(25, 13).Value = 1.059269350678369
